$wb = $excel.ActiveWorkbook

$aboutWs  = $wb.Worksheets.Item("About")
$fpiebpWs = $wb.Worksheets.Item("FPIEBP")

# --- "About" sheet: refresh the source/updated date shown in C1 ---
# (45294 -> 45379, i.e. 2024-01-03 -> 2024-03-28), and scroll the sheet
# down a bit so row 6 sits at the top when it is reopened.
$aboutWs.Activate()
$aboutWs.Range("C1").Value = 45379
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

# --- "FPIEBP" sheet: update the balancing priorities for "hard coal" ---
# production priority 3 -> 1, imports priority 2 -> 3, exports priority 1 -> 2
$fpiebpWs.Activate()
$fpiebpWs.Range("B3").Value = 1
$fpiebpWs.Range("C3").Value = 3
$fpiebpWs.Range("D3").Value = 2

# Leave the cursor on E3 (matches the saved selection) and keep FPIEBP
# as the active/displayed tab.
$fpiebpWs.Range("E3").Select() | Out-Null
